# Rerun all modules with HGG mis-assignments resolved
# Updates recomputed p-values in column C across all histology sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # Low-grade glioma
$ws.Range("C3").Value = 0.605339466053395
$ws.Range("C6").Value = 0.0717928207179282
$ws.Range("C7").Value = 0.22017798220178
$ws.Range("C8").Value = 0.783321667833217

$ws = $wb.Worksheets.Item(2)  # Non-neoplastic tumor
$ws.Range("C3").Value = 0.910708929107089
$ws.Range("C5").Value = 0.0000999900009999
$ws.Range("C6").Value = 0.257574242575742
$ws.Range("C7").Value = 0.857514248575142

$ws = $wb.Worksheets.Item(3)  # Mixed neuronal-glial tumor
$ws.Range("C3").Value = 0.2001799820018
$ws.Range("C6").Value = 0.526847315268473
$ws.Range("C7").Value = 0.794020597940206
$ws.Range("C8").Value = 0.695230476952305

$ws = $wb.Worksheets.Item(4)  # Medulloblastoma
$ws.Range("C3").Value = 0.777222277772223
$ws.Range("C6").Value = 0.46965303469653
$ws.Range("C7").Value = 0.293670632936706
$ws.Range("C8").Value = 0.83981601839816

$ws = $wb.Worksheets.Item(5)  # Schwannoma
$ws.Range("C3").Value = 0.805019498050195
$ws.Range("C5").Value = 0.505049495050495
$ws.Range("C6").Value = 0.931806819318068
$ws.Range("C7").Value = 0.848115188481152

$ws = $wb.Worksheets.Item(6)  # Mesenchymal tumor
$ws.Range("C3").Value = 0.127087291270873
$ws.Range("C5").Value = 0.0005999400059994
$ws.Range("C6").Value = 0.0216978302169783
$ws.Range("C7").Value = 0.844015598440156
$ws.Range("C8").Value = 0.48965103489651

$ws = $wb.Worksheets.Item(7)  # Germ cell tumor
$ws.Range("C3").Value = 0.615738426157384
$ws.Range("C4").Value = 0.0000999900009999
$ws.Range("C5").Value = 0.0106989301069893
$ws.Range("C6").Value = 0.856714328567143
$ws.Range("C7").Value = 0.359264073592641
$ws.Range("C8").Value = 0.275672432756724

$ws = $wb.Worksheets.Item(8)  # Craniopharyngioma
$ws.Range("C3").Value = 0.738426157384262
$ws.Range("C5").Value = 0.0063993600639936
$ws.Range("C6").Value = 0.163283671632837
$ws.Range("C7").Value = 0.573942605739426

$ws = $wb.Worksheets.Item(9)  # Other tumor
$ws.Range("C3").Value = 0.0477952204779522
$ws.Range("C5").Value = 0.0005999400059994
$ws.Range("C6").Value = 0.877212278772123
$ws.Range("C7").Value = 0.991300869913009

$ws = $wb.Worksheets.Item(10)  # Ependymoma
$ws.Range("C3").Value = 0.391960803919608
$ws.Range("C6").Value = 0.266273372662734
$ws.Range("C7").Value = 0.286271372862714
$ws.Range("C8").Value = 0.965803419658034

$ws = $wb.Worksheets.Item(11)  # DIPG or DMG
$ws.Range("C3").Value = 0.492150784921508
$ws.Range("C6").Value = 0.0148985101489851
$ws.Range("C7").Value = 0.622837716228377

$ws = $wb.Worksheets.Item(12)  # ATRT
$ws.Range("C3").Value = 0.377762223777622
$ws.Range("C5").Value = 0.0004999500049995
$ws.Range("C6").Value = 0.503349665033497
$ws.Range("C7").Value = 0.277272272772723
$ws.Range("C8").Value = 0.171482851714829

$ws = $wb.Worksheets.Item(13)  # Other high-grade glioma
$ws.Range("C3").Value = 0.797620237976202
$ws.Range("C6").Value = 0.831916808319168
$ws.Range("C7").Value = 0.0171982801719828
$ws.Range("C8").Value = 0.661533846615338

$ws = $wb.Worksheets.Item(14)  # Meningioma
$ws.Range("C3").Value = 0.941805819418058
$ws.Range("C5").Value = 0.0004999500049995
$ws.Range("C6").Value = 0.681231876812319
$ws.Range("C7").Value = 0.314868513148685

$ws = $wb.Worksheets.Item(15)  # Neurofibroma plexiform
$ws.Range("C5").Value = 0.153584641535846
$ws.Range("C6").Value = 0.0543945605439456
$ws.Range("C7").Value = 0.299270072992701

$ws = $wb.Worksheets.Item(16)  # Oligodendroglioma
$ws.Range("C3").Value = 0.174182581741826
$ws.Range("C5").Value = 0.204479552044796
$ws.Range("C7").Value = 0.211778822117788
